$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 14.36854692522047
$ws.Range("C2").Value = 12.38021239343275
$ws.Range("E2").Value = 13.06137553685492
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 21.85900228799304
$ws.Range("H2").Value = 12.72181295168823
$ws.Range("I2").Value = 19.81728697229302
$ws.Range("L2").Value = 9.947110547675999
$ws.Range("M2").Value = 13.92508164092805
$ws.Range("O2").Value = 18.41651148181071

$ws.Range("B3").Value = 13.70205755039805
$ws.Range("C3").Value = 12.20115111110407
$ws.Range("E3").Value = 13.12579840729361
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 22.06025434247227
$ws.Range("H3").Value = 12.78862197609126
$ws.Range("I3").Value = 19.96867562853138
$ws.Range("L3").Value = 9.955489349794986
$ws.Range("M3").Value = 13.76584757187563
$ws.Range("O3").Value = 18.54268540590802

$ws.Range("B4").Value = 13.27523495086323
$ws.Range("C4").Value = 12.09022848278583
$ws.Range("E4").Value = 13.16745884363122
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 22.19407283931687
$ws.Range("H4").Value = 12.83207369598199
$ws.Range("I4").Value = 20.0663903939381
$ws.Range("L4").Value = 9.962017909403553
$ws.Range("M4").Value = 13.66808599680153
$ws.Range("O4").Value = 18.6251803314928

$ws.Range("B5").Value = 13.09706544313829
$ws.Range("C5").Value = 12.04482343165117
$ws.Range("E5").Value = 13.18496652527727
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 22.25116498426289
$ws.Range("H5").Value = 12.85039244320006
$ws.Range("I5").Value = 20.10741033457606
$ws.Range("L5").Value = 9.965026882052753
$ws.Range("M5").Value = 13.62828462814532
$ws.Range("O5").Value = 18.66005966119846

$ws.Range("B6").Value = 13.06723079975119
$ws.Range("C6").Value = 12.03727295078154
$ws.Range("E6").Value = 13.18790576094012
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 22.26079921425432
$ws.Range("H6").Value = 12.85347122754252
$ws.Range("I6").Value = 20.11429425465962
$ws.Range("L6").Value = 9.965547584251212
$ws.Range("M6").Value = 13.62167891594408
$ws.Range("O6").Value = 18.66592753271707

$ws.Range("B7").Value = 13.27284897417387
$ws.Range("C7").Value = 12.08961689997533
$ws.Range("E7").Value = 13.16769280748447
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 22.19483246099328
$ws.Range("H7").Value = 12.83231827097622
$ws.Range("I7").Value = 20.06693873884012
$ws.Range("L7").Value = 9.962057077618978
$ws.Range("M7").Value = 13.66754902444616
$ws.Range("O7").Value = 18.62564561822728

$ws.Range("B8").Value = 14.14249758609642
$ws.Range("C8").Value = 12.31870055604495
$ws.Range("E8").Value = 13.08315265973703
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 21.9262564425008
$ws.Range("H8").Value = 12.74434456707947
$ws.Range("I8").Value = 19.86849937314585
$ws.Range("L8").Value = 9.949712656742651
$ws.Range("M8").Value = 13.87019693187512
$ws.Range("O8").Value = 18.45897282692384

$ws.Range("B9").Value = 15.7013771168245
$ws.Range("C9").Value = 12.75833330888829
$ws.Range("E9").Value = 12.93399889167358
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 21.4816784657137
$ws.Range("H9").Value = 12.59108765662678
$ws.Range("I9").Value = 19.51700003741524
$ws.Range("L9").Value = 9.936461202510106
$ws.Range("M9").Value = 14.26612503315019
$ws.Range("O9").Value = 18.17204569482694

$ws.Range("B10").Value = 16.74984920665653
$ws.Range("C10").Value = 13.07309377354762
$ws.Range("E10").Value = 12.83445583410334
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 21.20617807305492
$ws.Range("H10").Value = 12.4901924853882
$ws.Range("I10").Value = 19.28150066474069
$ws.Range("L10").Value = 9.933368002037749
$ws.Range("M10").Value = 14.55408810231274
$ws.Range("O10").Value = 17.98565644226765

$ws.Range("B11").Value = 17.20461920369613
$ws.Range("C11").Value = 13.21401735816459
$ws.Range("E11").Value = 12.79133054195387
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 21.09218612489348
$ws.Range("H11").Value = 12.44682615232396
$ws.Range("I11").Value = 19.17926310338223
$ws.Range("L11").Value = 9.933394027862951
$ws.Range("M11").Value = 14.68403278409715
$ws.Range("O11").Value = 17.90618445706001

$ws.Range("B12").Value = 17.37356197757233
$ws.Range("C12").Value = 13.26701789621204
$ws.Range("E12").Value = 12.77530878318376
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 21.05067036373958
$ws.Range("H12").Value = 12.43076799902094
$ws.Range("I12").Value = 19.14124881368615
$ws.Range("L12").Value = 9.933609043015972
$ws.Range("M12").Value = 14.73305384527028
$ws.Range("O12").Value = 17.87685710005552

$ws.Range("B13").Value = 17.33732361060203
$ws.Range("C13").Value = 13.25562007862961
$ws.Range("E13").Value = 12.7787456409519
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 21.05953778794106
$ws.Range("H13").Value = 12.43421023642046
$ws.Range("I13").Value = 19.14940474212661
$ws.Range("L13").Value = 9.933553626683846
$ws.Range("M13").Value = 14.72250511037303
$ws.Range("O13").Value = 17.88313912220236

$ws.Range("B14").Value = 17.21858413853714
$ws.Range("C14").Value = 13.21838525913146
$ws.Range("E14").Value = 12.79000624013887
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 21.08873740055674
$ws.Range("H14").Value = 12.4454977500334
$ws.Range("I14").Value = 19.17612161372095
$ws.Range("L14").Value = 9.933407611157961
$ws.Range("M14").Value = 14.68806970473832
$ws.Range("O14").Value = 17.90375629303524

$ws.Range("B15").Value = 17.14542500275232
$ws.Range("C15").Value = 13.19552930846517
$ws.Range("E15").Value = 12.7969438636628
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 21.106838568002
$ws.Range("H15").Value = 12.45245903960578
$ws.Range("I15").Value = 19.19257766321873
$ws.Range("L15").Value = 9.933344861503739
$ws.Range("M15").Value = 14.66695175599314
$ws.Range("O15").Value = 17.9164848570703

$ws.Range("B16").Value = 16.71967545478731
$ws.Range("C16").Value = 13.0638352021457
$ws.Range("E16").Value = 12.83731746725549
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 21.21385755276818
$ws.Range("H16").Value = 12.49307749099814
$ws.Range("I16").Value = 19.28828034775923
$ws.Range("L16").Value = 9.933395062452519
$ws.Range("M16").Value = 14.54557174632988
$ws.Range("O16").Value = 17.9909573191628

$ws.Range("B17").Value = 16.45274925063275
$ws.Range("C17").Value = 12.98243750659256
$ws.Range("E17").Value = 12.86263693429327
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 21.28242788490989
$ws.Range("H17").Value = 12.51864373841451
$ws.Range("I17").Value = 19.34824197909897
$ws.Range("L17").Value = 9.933792288286343
$ws.Range("M17").Value = 14.47081462434462
$ws.Range("O17").Value = 18.03800729519545

$ws.Range("B18").Value = 16.29713589192679
$ws.Range("C18").Value = 12.93540896503441
$ws.Range("E18").Value = 12.877403162638
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 21.32293334051664
$ws.Range("H18").Value = 12.53358704813293
$ws.Range("I18").Value = 19.38319096429685
$ws.Range("L18").Value = 9.934155716155161
$ws.Range("M18").Value = 14.42771940697929
$ws.Range("O18").Value = 18.06556955167222

$ws.Range("B19").Value = 16.24409250196172
$ws.Range("C19").Value = 12.9194509238383
$ws.Range("E19").Value = 12.88243768536938
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 21.33683018755799
$ws.Range("H19").Value = 12.53868752983961
$ws.Range("I19").Value = 19.39510328705124
$ws.Range("L19").Value = 9.934301971917026
$ws.Range("M19").Value = 14.41311254581421
$ws.Range("O19").Value = 18.07498752577028

$ws.Range("B20").Value = 16.4813803488975
$ws.Range("C20").Value = 12.99112452746531
$ws.Range("E20").Value = 12.85992061947945
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 21.2750180202762
$ws.Range("H20").Value = 12.51589750869158
$ws.Range("I20").Value = 19.34181130579212
$ws.Range("L20").Value = 9.933736042117596
$ws.Range("M20").Value = 14.47878293264105
$ws.Range("O20").Value = 18.03294694365222

$ws.Range("B21").Value = 17.25355006980678
$ws.Range("C21").Value = 13.22933221259774
$ws.Range("E21").Value = 12.78669035786268
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 21.0801158065853
$ws.Range("H21").Value = 12.44217246400853
$ws.Range("I21").Value = 19.16825522078704
$ws.Range("L21").Value = 9.933444939240248
$ws.Range("M21").Value = 14.69818954521728
$ws.Range("O21").Value = 17.8976796960666

$ws.Range("B22").Value = 17.73912779911684
$ws.Range("C22").Value = 13.38287526341328
$ws.Range("E22").Value = 12.74062975551217
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 20.96236612644912
$ws.Range("H22").Value = 12.39610891256663
$ws.Range("I22").Value = 19.05891021079787
$ws.Range("L22").Value = 9.934450029494046
$ws.Range("M22").Value = 14.84048373544322
$ws.Range("O22").Value = 17.81374631553919

$ws.Range("B23").Value = 17.48173407220829
$ws.Range("C23").Value = 13.30113459467253
$ws.Range("E23").Value = 12.76504895122686
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 21.02432331675289
$ws.Range("H23").Value = 12.42049999634125
$ws.Range("I23").Value = 19.11689688576839
$ws.Range("L23").Value = 9.933804545854818
$ws.Range("M23").Value = 14.76465069045364
$ws.Range("O23").Value = 17.85813319255319

$ws.Range("B24").Value = 16.46844293228846
$ws.Range("C24").Value = 12.98719784475981
$ws.Range("E24").Value = 12.86114801201652
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 21.27836464674904
$ws.Range("H24").Value = 12.51713831617128
$ws.Range("I24").Value = 19.34471712951294
$ws.Range("L24").Value = 9.933761050253779
$ws.Range("M24").Value = 14.47518082019387
$ws.Range("O24").Value = 18.03523313204316

$ws.Range("B25").Value = 15.29620385790961
$ws.Range("C25").Value = 12.64067970669877
$ws.Range("E25").Value = 12.9725787787003
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 21.59304879467273
$ws.Range("H25").Value = 12.63048994809079
$ws.Range("I25").Value = 19.60808104155602
$ws.Range("L25").Value = 9.938876623943942
$ws.Range("M25").Value = 14.1593842038862
$ws.Range("O25").Value = 18.17204569482694
